$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row shared strings -------------------------------------
# "<Name>_old"  -> "<Name>_FV2404"
# "<Name>_new"  -> "<Name>_FV2410"
# The header row (row 1) already carries direct cell formatting (bold,
# centered, shaded, bordered). Save a copy of that formatting in a scratch
# area first, so it can be restored bit-for-bit after the table is added
# below -- ListObjects.Add() otherwise "captures" the header's current
# format into a brand-new dxf/table-style, which is not part of the target
# change.
$hdr = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")
$hdr.Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$hdr.Style = "Normal"

$ws.Range("A1").Value = "Segmentname_FV2404"
$ws.Range("B1").Value = "Segmentgruppe_FV2404"
$ws.Range("C1").Value = "Segment_FV2404"
$ws.Range("D1").Value = "Datenelement_FV2404"
$ws.Range("E1").Value = "Segment ID_FV2404"
$ws.Range("F1").Value = "Code_FV2404"
$ws.Range("G1").Value = "Qualifier_FV2404"
$ws.Range("H1").Value = "Beschreibung_FV2404"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value = "Bedingung_FV2404"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2410"
$ws.Range("M1").Value = "Segmentgruppe_FV2410"
$ws.Range("N1").Value = "Segment_FV2410"
$ws.Range("O1").Value = "Datenelement_FV2410"
$ws.Range("P1").Value = "Segment ID_FV2410"
$ws.Range("Q1").Value = "Code_FV2410"
$ws.Range("R1").Value = "Qualifier_FV2410"
$ws.Range("S1").Value = "Beschreibung_FV2410"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value = "Bedingung_FV2410"

# --- Turn the data range into a real table (Table1) -----------------------
$range = $ws.Range("A1:U85")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Restore the header row's original direct formatting.
$scratch.Copy() | Out-Null
$hdr.PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$scratch.Clear() | Out-Null

# --- Freeze the header row -------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
